# Slide 3 ("Contents") - Content Placeholder 2:
#   "Best Practices"   -> split into "Best " / "Practices" runs,
#                          then a new "Resources" paragraph is added after it.
#   "Where to Start?"  -> split into "Where " / "to Start?" runs,
#                          then a new "Questions" paragraph is added after it.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$tr = $s.Shapes.Item(2).TextFrame.TextRange

# Split "Best Practices" into two runs: "Best " + "Practices"
$tr.Paragraphs(2).Characters(1, 5).Text = "Best "

# Split "Where to Start?" into two runs: "Where " + "to Start?"
$tr.Paragraphs(3).Characters(1, 6).Text = "Where "

# Insert a new paragraph "Resources" right after the "Best Practices" paragraph
[void]$tr.Paragraphs(2).InsertAfter([char]13 + "Resources")

# Insert a new paragraph "Questions" right after the "Where to Start?" paragraph
# (now paragraph 4, since "Resources" shifted everything below it down by one)
[void]$tr.Paragraphs(4).InsertAfter([char]13 + "Questions")
